# Update odds values on row 3 of Sheet1 to reflect the latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1.45
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 7.5
$ws.Range("J3").Value = 2.05
$ws.Range("K3").Value = 2.2
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8

$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 5.5
$ws.Range("Z3").Value = 9

$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 26

$ws.Range("AL3").Value = 67

$ws.Range("AW3").Value = 8.5
$ws.Range("AY3").Value = 51
